$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 20 of data
$ws.Range("A2").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$ws.Range("A20").Value = 45986

$ws.Range("B20").Value = 2025
$ws.Range("C20").Value = -2.451276118722334
$ws.Range("D20").Value = 2026
$ws.Range("E20").Value = -0.8888225292121632
